$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet immediately before "2022-Q3" and fill it
#    with the quarterly fund-holdings table.
# ---------------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($sheetQ3)
$newSheet.Name = "2022-Q4"

# Header row
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: row, A(index), B(code), C(name), D(scale), E(position), F(ratio), G(value), H(rank)
$fundRows = @(
    @("2","0","002685","中欧丰泓沪港深灵活配置混合A","54.86","93.29","6.48","3.5549","3"),
    @("3","1","005847","富国沪港深业绩驱动混合A","36.59","88.31","8.29","3.0333","3"),
    @("4","2","501087","交银施罗德瑞丰混合（LOF）","21.14","88.79","7.48","1.5813","3"),
    @("5","3","002686","中欧丰泓沪港深灵活配置混合C","16.87","93.29","6.48","1.0932","3"),
    @("6","4","010583","富国蓝筹精选股票（QDII）美元","13.38","91.97","7.85","1.0503","1"),
    @("7","5","007455","富国蓝筹精选股票（QDII）人民币","13.38","91.97","7.85","1.0503","1"),
    @("8","6","016464","兴证全球合瑞混合A","26.69","85.54","2.97","0.7927","8"),
    @("9","7","001605","国富沪港深成长精选股票","21.45","90.72","3.68","0.7894","3"),
    @("10","8","000934","国富大中华精选混合（QDII）","20.61","87.91","3.47","0.7152","7"),
    @("11","9","006370","国富大中华精选混合（QDII）美元","20.61","87.91","3.47","0.7152","7"),
    @("12","10","009846","富兰克林国海港股通远见价值混合","15.17","90.16","3.72","0.5643","3"),
    @("13","11","011117","富国沪港深业绩驱动混合C","6.15","88.31","8.29","0.5098","3"),
    @("14","12","016465","兴证全球合瑞混合C","13.89","85.54","2.97","0.4125","8"),
    @("15","13","013991","中欧港股通精选一年持有混合A","7.23","91.98","5.38","0.3890","6"),
    @("16","14","012744","光大保德信品质生活混合A","6.13","84.35","6.23","0.3819","2"),
    @("17","15","010088","工银优质成长混合A","15.36","77.14","2.23","0.3425","8"),
    @("18","16","013992","中欧港股通精选一年持有混合C","4.91","91.98","5.38","0.2642","6"),
    @("19","17","011635","富国港股通策略精选混合A","7.15","85.60","2.86","0.2045","6"),
    @("20","18","012758","光大保德信品质生活混合C","2.19","84.35","6.23","0.1364","2"),
    @("21","19","012770","光大保德信创新生活混合","2.89","88.40","4.55","0.1315","4"),
    @("22","20","006039","国富估值优势混合A","4.06","83.32","2.86","0.1161","10"),
    @("23","21","012060","富国全球消费精选混合（QDII）A","3.09","87.16","3.01","0.0930","5"),
    @("24","22","012061","富国全球消费精选混合（QDII）美元现汇","3.09","87.16","3.01","0.0930","5"),
    @("25","23","016298","中欧丰泰港股通混合C","0.43","93.01","6.25","0.0269","3"),
    @("26","24","011636","富国港股通策略精选混合C","0.73","85.60","2.86","0.0209","6"),
    @("27","25","010089","工银优质成长混合C","0.90","77.14","2.23","0.0201","8"),
    @("28","26","016297","中欧丰泰港股通混合A","0.27","93.01","6.25","0.0169","3"),
    @("29","27","001942","前海开源沪港深汇鑫灵活配置混合A","0.31","90.30","4.92","0.0153","6"),
    @("30","28","001943","前海开源沪港深汇鑫灵活配置混合C","0.27","90.30","4.92","0.0133","6"),
    @("31","29","014214","光大保德信核心资产混合A","0.32","85.79","3.57","0.0114","7"),
    @("32","30","000761","国富健康优质生活股票","0.11","86.45","3.38","0.0037","7"),
    @("33","31","014215","光大保德信核心资产混合C","0.03","85.79","3.57","0.0011","7"),
    @("34","32","017451","国富估值优势混合C","0.00","83.32","2.86","0","10"),
    @("35","33","012062","富国全球消费精选混合（QDII）C","-2.96","87.16","3.01","-0.0891","5")
)

foreach ($r in $fundRows) {
    $rowNum = [int]$r[0]

    $cellA = $newSheet.Cells.Item($rowNum, 1)
    $cellA.Value = [int]$r[1]
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1

    $cellB = $newSheet.Cells.Item($rowNum, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r[2]

    $cellC = $newSheet.Cells.Item($rowNum, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $r[3]

    $cellD = $newSheet.Cells.Item($rowNum, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $r[4]

    $cellE = $newSheet.Cells.Item($rowNum, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $r[5]

    $cellF = $newSheet.Cells.Item($rowNum, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $r[6]

    $cellG = $newSheet.Cells.Item($rowNum, 7)
    if ($r[7] -eq "0") {
        $cellG.Value = 0
    } else {
        $cellG.NumberFormat = "@"
        $cellG.Value = $r[7]
    }

    $cellH = $newSheet.Cells.Item($rowNum, 8)
    $cellH.Value = [int]$r[8]
}

# ---------------------------------------------------------------------------
# 2. Insert the new 2022-Q4 summary row at the top of the "总计" sheet and
#    reflow the existing quarters beneath it.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")
$ws1.Rows.Item(2).Insert()

$totalsRows = @(
    @(2,0,"2022-Q4",34,18.06),
    @(3,1,"2022-Q3",33,16.26),
    @(4,2,"2022-Q2",19,18.62),
    @(5,3,"2022-Q1",22,20.3),
    @(6,4,"2021-Q4",12,13.44),
    @(7,5,"2021-Q3",10,13.53),
    @(8,6,"2021-Q2",20,14.32),
    @(9,7,"2021-Q1",16,6.89)
)

foreach ($r in $totalsRows) {
    $rowNum = $r[0]

    $cellA = $ws1.Cells.Item($rowNum, 1)
    $cellA.Value = $r[1]
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1

    $ws1.Cells.Item($rowNum, 2).NumberFormat = "General"
    $ws1.Cells.Item($rowNum, 2).Value = $r[2]
    $ws1.Cells.Item($rowNum, 3).NumberFormat = "General"
    $ws1.Cells.Item($rowNum, 3).Value = $r[3]
    $ws1.Cells.Item($rowNum, 4).NumberFormat = "General"
    $ws1.Cells.Item($rowNum, 4).Value = $r[4]
}
